$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Program Capacity" values for IT1..IT10 (column B, rows 2..11).
# These are numeric-looking strings that must remain stored as text
# (shared strings), matching the original workbook's cell type.
$values = @{
    "B2"  = "89"   # IT1
    "B3"  = "45"   # IT2
    "B4"  = "27"   # IT3
    "B5"  = "44"   # IT4
    "B6"  = "85"   # IT5
    "B7"  = "6"    # IT6
    "B8"  = "13"   # IT7
    "B9"  = "29"   # IT8
    "B10" = "47"   # IT9
    "B11" = "44"   # IT10
}

foreach ($addr in $values.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so the numeric-looking value is kept as a string
    # (shared string) rather than being coerced into a numeric cell.
    $cell.NumberFormat = "@"
    $cell.Value = $values[$addr]
}

# Remove the temporary text number-format so the cells keep using the
# workbook's default style (no new style index left referenced on cells).
$ws.Range("B2:B11").ClearFormats()
